$d = $word.ActiveDocument

# --- Part 1: merge the bookmark-split run back into a single run ---
# Remove the now-redundant bookmark text split by replacing the two
# adjacent fragments (with the bookmark in between) with one continuous run.
$d.Content.Find.Execute(
    " much easier and cleaner to write A tab strip",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " much easier and cleaner to write A tab strip with tag helpers. however  writing kendo buttons with html helpers is easier to pass parameters through for some JavaScript action ",
    2) | Out-Null

# The old trailing text (that used to follow the bookmark) is now
# duplicated immediately after our replacement; remove it.
$d.Content.Find.Execute(
    "JavaScript action  with tag helpers. however  writing kendo buttons with html helpers is easier to pass parameters through for some JavaScript action ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "JavaScript action ",
    2) | Out-Null

# Remove the old bookmark (it will be re-created further up in the doc)
foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

# --- Part 2: append "." to the "Partial cshtml scripts..." paragraph and add a new paragraph after it ---
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Partial cshtml scripts*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $endRange = $target.Range
    $endRange.SetRange($endRange.End - 1, $endRange.End - 1)
    $endRange.InsertAfter(".")

    $paraEnd = $target.Range
    $paraEnd.SetRange($target.Range.End, $target.Range.End)
    $paraEnd.InsertParagraphAfter()

    $newPara = $target.Next()
    $newPara.Range.InsertAfter("This goes for templates and template Ids in grids etc. if there is any styling on grids/partials, it must be done either in-line or in style tags")

    $newEnd = $newPara.Range
    $newEnd.SetRange($newEnd.End - 1, $newEnd.End - 1)
    $d.Bookmarks.Add("_GoBack", $newEnd)
}
